$wb = $excel.ActiveWorkbook

# The "General Settings" sheet had a column ("LNE File Cabinet Path") removed.
$ws = $wb.Worksheets.Item("General Settings")
$ws.Activate()

# Select the entire column I (as the author did before deleting it) and delete it,
# shifting the remaining columns (Download folder name, EP Installation timeout,
# From EP service start until logs show EP active timeout, EP Service Timeout) left.
$ws.Columns("I").Select()
$ws.Columns("I").Delete()

# Leave the selection on the former column I (now holding what used to be column J).
$ws.Range("I1").Select()
